# Omni_ERP.docx cleanup
# ----------------------
# Removes, throughout the whole document body:
#   1. The 3 inline images (the screenshot drawings).
#   2. The 12 "────…" separator-rule paragraphs
#      (pPr: spacing before=120/after=120, single run, color CCCCCC, sz 16).
#   3. The 10 empty spacer paragraphs that sit right after each code table
#      (pPr: spacing before=40, no runs at all).
#
# Matching is done structurally (inline-shape presence, exact separator
# text, "empty + SpaceBefore=40(=2pt)") rather than by hard-coded paragraph
# index, since removing each paragraph shifts every later index. We first
# collect every paragraph index that has to go, then delete from the end
# of the document backwards so earlier indices stay valid while we work.

$d = $word.ActiveDocument

$separatorText = "────────────────────────────────────────────────────────────"

$count = $d.Paragraphs.Count
$indicesToDelete = @()

for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range

    # Never touch anything living inside a table cell.
    if ($rng.Information(12)) {
        continue
    }

    $text = $rng.Text
    $trimmedText = $text.TrimEnd([char]13, [char]7)

    $hasInlineImage = ($rng.InlineShapes.Count -gt 0)
    $isSeparatorRule = ($trimmedText -eq $separatorText)
    $isEmptySpacer = (($trimmedText.Length -eq 0) -and ($para.Format.SpaceBefore -eq 2))

    if ($hasInlineImage -or $isSeparatorRule -or $isEmptySpacer) {
        $indicesToDelete += $i
    }
}

Write-Output "Paragraphs scheduled for removal: $($indicesToDelete.Count)"

for ($k = $indicesToDelete.Count - 1; $k -ge 0; $k--) {
    $idx = $indicesToDelete[$k]
    $d.Paragraphs($idx).Range.Delete()
}

Write-Output "Remaining paragraph count: $($d.Paragraphs.Count)"
